# Insert two new weekly price-report rows for "Acelga" (Vega Central Mapocho de
# Santiago) just before the existing row that starts the current 772 block.
# This shifts the old rows 772:806 down to 774:808 (matches new dimension
# A1:R808) and fills the two freshly-inserted rows (772, 773) with the new
# week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 772/773, pushing everything from 772 downward by two.
$ws.Range("A772:A773").EntireRow.Insert()

# New row 772: Acelga, Primera, market date 2023-03-23 (serial 45008)
$ws.Cells.Item(772,1).Value = 9
$ws.Cells.Item(772,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(772,3).Value = "Metropolitana"
$ws.Cells.Item(772,4).Value = 45008
$ws.Cells.Item(772,5).Value = 13
$ws.Cells.Item(772,6).Value = 100112009
$ws.Cells.Item(772,7).Value = "Acelga"
$ws.Cells.Item(772,8).Value = "Sin especificar"
$ws.Cells.Item(772,9).Value = "Primera"
$ws.Cells.Item(772,10).Value = 70
$ws.Cells.Item(772,11).Value = 16000
$ws.Cells.Item(772,12).Value = 16000
$ws.Cells.Item(772,13).Value = 16000
$ws.Cells.Item(772,14).Value = "`$/docena de atados"
$ws.Cells.Item(772,15).Value = "Región Metropolitana"
$ws.Cells.Item(772,16).Value = 5333
$ws.Cells.Item(772,17).Value = 3
$ws.Cells.Item(772,18).Value = "Hortaliza"

# New row 773: Acelga, Segunda, market date 2023-03-23 (serial 45008)
$ws.Cells.Item(773,1).Value = 9
$ws.Cells.Item(773,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(773,3).Value = "Metropolitana"
$ws.Cells.Item(773,4).Value = 45008
$ws.Cells.Item(773,5).Value = 13
$ws.Cells.Item(773,6).Value = 100112009
$ws.Cells.Item(773,7).Value = "Acelga"
$ws.Cells.Item(773,8).Value = "Sin especificar"
$ws.Cells.Item(773,9).Value = "Segunda"
$ws.Cells.Item(773,10).Value = 34
$ws.Cells.Item(773,11).Value = 13000
$ws.Cells.Item(773,12).Value = 13000
$ws.Cells.Item(773,13).Value = 13000
$ws.Cells.Item(773,14).Value = "`$/docena de atados"
$ws.Cells.Item(773,15).Value = "Región Metropolitana"
$ws.Cells.Item(773,16).Value = 4333
$ws.Cells.Item(773,17).Value = 3
$ws.Cells.Item(773,18).Value = "Hortaliza"
